$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.676.86"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").Value = "1.852.77"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4262"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.86%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3648"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.50%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07302"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8750"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.66%  "

$ws.Range("D13").Value = "1.832.97"
$ws.Range("E13").Value = "  +1.06%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.316"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.03%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.516"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06890"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009017"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").Value = "27.685.73"
$ws.Range("E22").Value = "  +0.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.974"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.61%  "

$ws.Range("E24").Value = "  -2.14%  "

$ws.Range("D25").Value = "2.080.75"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.984"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.84%  "

$ws.Range("E28").Value = "  +3.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "122.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.266"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.864"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08913"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7685"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.03%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.966"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.79%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.525"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.106"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05384"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01937"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.826"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5061"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.805"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1647"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.361"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06537"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "104.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4665"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("E50").Value = "  -0.89%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.33%  "
